$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed Price (D) / Volume(1h) (E) figures from the latest crawl.
# These columns hold plain text (e.g. "303.68", "4.27%"), so each cell is
# switched to Text format before the new literal is written (stopping Excel
# from reinterpreting it as a number/percentage) and then restored to the
# workbook default "Normal" style, leaving formatting untouched.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "303.68"
Set-TextValue "E2" "4.27%"
Set-TextValue "D3" "32.08"
Set-TextValue "E3" "9.15%"
Set-TextValue "D4" "5.263"
Set-TextValue "E4" "0.20%"
Set-TextValue "D5" "0.07524"
Set-TextValue "E5" "5.28%"
Set-TextValue "D6" "7.913"
Set-TextValue "E6" "5.69%"
Set-TextValue "E7" "6.71%"
Set-TextValue "D8" "1.498"
Set-TextValue "E8" "6.85%"
Set-TextValue "D9" "0.9220"
Set-TextValue "E9" "1.37%"
Set-TextValue "D10" "0.1695"
Set-TextValue "E10" "4.42%"
Set-TextValue "D11" "0.08004"
Set-TextValue "E11" "5.14%"
Set-TextValue "D12" "0.08008"
Set-TextValue "E12" "3.35%"
Set-TextValue "D13" "0.03036"
Set-TextValue "E13" "4.00%"
Set-TextValue "D14" "0.09914"
Set-TextValue "E14" "10.30%"
Set-TextValue "D15" "0.001520"
Set-TextValue "E15" "-4.52%"
Set-TextValue "D16" "0.04596"
Set-TextValue "E16" "1.72%"
Set-TextValue "D17" "0.006538"
Set-TextValue "E17" "3.23%"
Set-TextValue "D18" "3.462"
Set-TextValue "E18" "-1.26%"
Set-TextValue "D19" "2.227"
Set-TextValue "E19" "-0.14%"
Set-TextValue "D20" "0.3300"
Set-TextValue "E20" "1.12%"
Set-TextValue "D21" "0.1336"
Set-TextValue "E21" "-1.51%"
Set-TextValue "D22" "4.475"
Set-TextValue "E22" "11.75%"
Set-TextValue "E23" "1.63%"
Set-TextValue "D24" "0.001216"
Set-TextValue "E24" "1.29%"
Set-TextValue "D25" "0.004456"
Set-TextValue "E25" "5.23%"
Set-TextValue "D26" "0.0001397"
Set-TextValue "E26" "20.58%"
Set-TextValue "D27" "0.0001773"
Set-TextValue "E27" "-7.58%"
Set-TextValue "D39" "0.01711"
Set-TextValue "E39" "2,534.07%"
Set-TextValue "D40" "0.04491"
Set-TextValue "E40" "2.35%"
Set-TextValue "D41" "0.006985"
Set-TextValue "E41" "0.13%"
Set-TextValue "D42" "0.1353"
Set-TextValue "E42" "6.79%"
Set-TextValue "D43" "0.002136"
Set-TextValue "E43" "1.75%"
Set-TextValue "D44" "0.01289"
Set-TextValue "E44" "-2.73%"
Set-TextValue "D45" "0.00006165"
Set-TextValue "E45" "6.06%"
Set-TextValue "D46" "1.863"
Set-TextValue "E46" "-3.43%"
Set-TextValue "E47" "16.28%"
